# fpspreadsheet test-data.xlsx: add an x=0 data point as the new first row of
# the sin/cos series on Sheet1 and Sheet2 (so charts built off this source can
# show a leading label/series-start at x=0), and tighten the display format of
# the computed (B) columns.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws3 = $wb.Worksheets.Item("Sheet3")

# --- Sheet1 (sin) --------------------------------------------------------
# Insert a new row 2 holding x=0, push the rest of the series down one row.
$ws1.Rows.Item(2).Insert()
$ws1.Range("A2:B2").ClearFormats()
$ws1.Range("A2").Value = 0
$ws1.Range("B2").Formula = "=SIN(A2*0.3)"

# New leading x cell is right aligned; the whole B (computed) column gets a
# tighter 5-decimal display format.
$ws1.Range("A2").HorizontalAlignment = -4152
$ws1.Range("B2:B32").NumberFormat = "0.00000"

# --- Sheet2 (cos) ---------------------------------------------------------
# Same pattern: insert x=0 as the new row 2.
$ws2.Rows.Item(2).Insert()
$ws2.Range("A2:B2").ClearFormats()
$ws2.Range("A2").Value = 0
$ws2.Range("B2").Formula = "=COS(A2*0.3)"

# Computed column gets a 3-decimal display format.
$ws2.Range("B2:B27").NumberFormat = "0.000"

# --- Selection bookkeeping -------------------------------------------------
# Reset the stale selections left over on the inactive sheets, then leave the
# cursor on the newly-added row on the active sheet (Sheet1).
[void]$ws2.Range("A1").Select()
[void]$ws3.Range("A1").Select()
$ws1.Activate()
[void]$ws1.Range("A2").Select()
